# Re-theme the deck: the "Integral" theme that was driving the slides is
# swapped back out for the stock "Office Theme" palette, and the one table
# on slide 5 is switched from the custom Table_0 style to a different
# built-in table style.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 (the single table in the deck) ---------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{3194CF84-A136-4752-8DC6-8613F5ADDD8D}")

# --- 2. Theme colour scheme: restore the plain "Office" colours -------
# (the deck's slide master / theme used the "Integral" - Red Violet -
# palette; swap every slot back to the standard Office theme colours)
$colorScheme = $slide.ColorScheme

$colorScheme.Colors(1).RGB  = 0         # dk1      000000
$colorScheme.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$colorScheme.Colors(3).RGB  = 6968388   # dk2      44546A
$colorScheme.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$colorScheme.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$colorScheme.Colors(6).RGB  = 3243501   # accent2  ED7D31
$colorScheme.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$colorScheme.Colors(8).RGB  = 49407     # accent4  FFC000
$colorScheme.Colors(9).RGB  = 12874308  # accent5  4472C4
$colorScheme.Colors(10).RGB = 4697456   # accent6  70AD47
$colorScheme.Colors(11).RGB = 12673797  # hlink    0563C1
$colorScheme.Colors(12).RGB = 7491477   # folHlink 954F72
